# business questions code is updated:
# Add a new test-case row (TC004) to Sheet1, mirroring the row above it
# (TC003) except for the new test-case id and the AA column answer
# (No -> Yes), wire up the same mailto hyperlink used by the other rows,
# and leave Sheet1 as the active/selected sheet with the new row's last
# cell selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# 1) Copy formatting from row 4 down onto row 5 first, so every cell
#    already carries the right number format / font / etc. before any
#    values are written (this keeps numeric-looking text such as zip
#    codes and monetary amounts stored as text instead of numbers).
$ws1.Range("A4:AC4").Copy()
$ws1.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Populate row 5 values - identical to row 4 except A5 (new test case
#    id) and AA5 (answer flips from No to Yes).
$ws1.Range("A5").Value = "TC004"
$ws1.Range("B5").Value = $ws1.Range("B4").Value2
$ws1.Range("C5").Value = $ws1.Range("C4").Value2
$ws1.Range("D5").Value = $ws1.Range("D4").Value2
$ws1.Range("E5").Value = $ws1.Range("E4").Value2
$ws1.Range("F5").Value = $ws1.Range("F4").Value2
$ws1.Range("G5").Value = $ws1.Range("G4").Text
$ws1.Range("H5").Value = $ws1.Range("H4").Value2
$ws1.Range("I5").Value = $ws1.Range("I4").Value2
$ws1.Range("J5").Value = $ws1.Range("J4").Value2
$ws1.Range("K5").Value = $ws1.Range("K4").Value2
$ws1.Range("L5").Value = $ws1.Range("L4").Value2
$ws1.Range("M5").Value = $ws1.Range("M4").Value2
$ws1.Range("N5").Value = $ws1.Range("N4").Value2
$ws1.Range("O5").Value = $ws1.Range("O4").Value2
$ws1.Range("P5").Value = $ws1.Range("P4").Value2
$ws1.Range("Q5").Value = $ws1.Range("Q4").Value2
$ws1.Range("R5").Value = $ws1.Range("R4").Value2
$ws1.Range("S5").Value = $ws1.Range("S4").Value2
$ws1.Range("T5").Value = $ws1.Range("T4").Text
$ws1.Range("U5").Value = $ws1.Range("U4").Value2
$ws1.Range("V5").Value = $ws1.Range("V4").Text
$ws1.Range("W5").Value = $ws1.Range("W4").Value2
$ws1.Range("X5").Value = $ws1.Range("X4").Value2
$ws1.Range("Y5").Value = $ws1.Range("Y4").Text
$ws1.Range("Z5").Value = $ws1.Range("Z4").Text
$ws1.Range("AA5").Value = "Yes"
$ws1.Range("AB5").Value = $ws1.Range("AB4").Value2
$ws1.Range("AC5").Value = $ws1.Range("AC4").Value2

# 3) Give the new email cell the same mailto hyperlink the other rows
#    have.
$ws1.Hyperlinks.Add($ws1.Range("B5"), "mailto:Email3@gmail.com")

# 4) Adding the hyperlink resets the cell's style, so re-apply the
#    Hyperlink format from the cell above (reuses the existing style).
$ws1.Range("B4").Copy()
$ws1.Range("B5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 5) Leave Sheet1 as the active/selected sheet (it had lost that status
#    to Sheet2), with the new row's last cell selected.
$ws2.Activate()
$ws1.Activate()
$ws1.Range("AA5").Select()
